$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helpers
# ---------------------------------------------------------------------------

function Get-ParaIndexByText($text) {
    $i = 0
    foreach ($p in $d.Paragraphs) {
        $i = $i + 1
        $t = $p.Range.Text
        $t = $t.TrimEnd([char]13, [char]7)
        if ($t -eq $text) {
            return $i
        }
    }
    return -1
}

function Set-ParagraphXml($paraIndex, $pPrXml, $runsXml) {
    $p = $d.Paragraphs($paraIndex)
    $start = $p.Range.Start
    $end = $p.Range.End
    $rr = $d.Range($start, $end)
    $xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
        '<w:body><w:p>' + $pPrXml + $runsXml + '</w:p></w:body>' +
        '</w:document></pkg:xmlData></pkg:part></pkg:package>'
    $rr.InsertXML($xml)
}

$pPr0 = '<w:pPr><w:pStyle w:val="Paragraphedeliste"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr>'
$pPr1 = '<w:pPr><w:pStyle w:val="Paragraphedeliste"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr>'

# ---------------------------------------------------------------------------
# 1) "Irving Berlin" question: rewrite the question stem and its three
#    answer choices (resolve indices up front, before any edits happen).
# ---------------------------------------------------------------------------

$idxBerlinQ      = Get-ParaIndexByText("Irving Berlin, author of the song God Bless America, also wrote the favorite Christmas song:")
$idxBerlinA1     = Get-ParaIndexByText("Winter Wonderland")
$idxBerlinA2     = Get-ParaIndexByText("Let It Snow")
$idxBerlinA3     = Get-ParaIndexByText("White Christmas")

# Answer 3: "The co-writer of the Broadway show The Sound of Music"
$runsA3 = '<w:r><w:t xml:space="preserve">The co-writer of the Broadway show </w:t></w:r>' +
          '<w:r><w:rPr><w:i/></w:rPr><w:t>The Sound of Music</w:t></w:r>'
Set-ParagraphXml $idxBerlinA3 $pPr1 $runsA3

# Answer 2: "The author of God Bless America and White Christmas"
$runsA2 = '<w:r><w:t xml:space="preserve">The author of </w:t></w:r>' +
          '<w:r><w:rPr><w:i/></w:rPr><w:t>God Bless America</w:t></w:r>' +
          '<w:r><w:t xml:space="preserve"> and </w:t></w:r>' +
          '<w:r><w:rPr><w:i/></w:rPr><w:t>White Christmas</w:t></w:r>'
Set-ParagraphXml $idxBerlinA2 $pPr1 $runsA2

# Answer 1: "A famous American opera singer"
$runsA1 = '<w:r><w:t>A famous American opera singer</w:t></w:r>'
Set-ParagraphXml $idxBerlinA1 $pPr1 $runsA1

# Question stem: "Who was Irving Berlin?"
$runsQ = '<w:r><w:t>Who was Irving Berlin?</w:t></w:r>'
Set-ParagraphXml $idxBerlinQ $pPr0 $runsQ

# ---------------------------------------------------------------------------
# 2) "Which one of the following is NOT a US state?" answer choices
# ---------------------------------------------------------------------------

$d.Content.Find.Execute("North Carolina", $true, $false, $false, $false, $false, $true, 1, $false, "New Hampshire", 2)
$d.Content.Find.Execute("Georgia", $true, $false, $false, $false, $false, $true, 1, $false, "New York", 2)

# ---------------------------------------------------------------------------
# 3) "second most spoken language" question stem
# ---------------------------------------------------------------------------

$idxSpoken = Get-ParaIndexByText("When the United States were founded, what was the second most spoken language?")
$runsSpoken = '<w:r><w:t>W</w:t></w:r>' +
              '<w:r><w:t>hat</w:t></w:r>' +
              '<w:r><w:t xml:space="preserve"> is</w:t></w:r>' +
              '<w:r><w:t xml:space="preserve"> the second most spoken language</w:t></w:r>' +
              '<w:r><w:t xml:space="preserve"> in the US</w:t></w:r>' +
              '<w:r><w:t>?</w:t></w:r>'
Set-ParagraphXml $idxSpoken $pPr0 $runsSpoken
